{"js": "// Remove the \"Cumulative GPA 3.1\" / \"13\" paragraph from the Education section,\n// and remove the \"Is a Next.js project.\" paragraph plus the blank list\n// paragraph that immediately follows it in the Project Experience section.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Collect the paragraphs to delete by matching on their trimmed text. We\n// gather them first (so index shifts caused by earlier deletes don't affect\n// later lookups) and then delete them in reverse document order.\nconst toDelete = [];\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (text === \"Cumulative GPA 3.113\" || text === \"Is a Next.js project.\") {\n    toDelete.push(i);\n    // The empty \"ListParagraph\" paragraph directly following the Next.js\n    // project bullet is also removed as part of this edit.\n    if (text === \"Is a Next.js project.\" && i + 1 < items.length && items[i + 1].text.trim() === \"\") {\n      toDelete.push(i + 1);\n    }\n  }\n}\n\ntoDelete.sort((a, b) => b - a);\nfor (const idx of toDelete) {\n  items[idx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Cumulative GPA 3.1\" / \"13\" paragraph from the Education section,\n# and remove the \"Is a Next.js project.\" paragraph plus the blank list\n# paragraph that immediately follows it in the Project Experience section.\n\n$d = $word.ActiveDocument\n\n# Collect paragraph indices (1-based, Word COM convention) to delete, based\n# on their trimmed text, walking back-to-front so earlier deletions don't\n# invalidate later indices.\n$count = $d.Paragraphs.Count\n$toDelete = @()\n\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($text -eq \"Cumulative GPA 3.113\" -or $text -eq \"Is a Next.js project.\") {\n        $toDelete += $i\n        # The empty \"ListParagraph\" paragraph directly following the Next.js\n        # project bullet is also removed as part of this edit.\n        if ($text -eq \"Is a Next.js project.\" -and ($i + 1) -le $count) {\n            $nextText = $d.Paragraphs.Item($i + 1).Range.Text.Trim()\n            if ($nextText -eq \"\") {\n                $toDelete += ($i + 1)\n            }\n        }\n    }\n}\n\n$toDelete = $toDelete | Sort-Object -Descending\n\nforeach ($idx in $toDelete) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
